$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the header row (shared strings), keep A1 ("id") as-is
$ws.Range("B1").Value = "first_name"
$ws.Range("C1").Value = "last_name"
$ws.Range("D1").Value = "date (dd-MM-yyyy)"
$ws.Range("E1").Value = "punchIn (HH:mm)"
$ws.Range("F1").Value = "punchOut (HH:mm)"

# Widen the data columns so the (now longer) headers are readable
$ws.Columns.Item(2).ColumnWidth = 14.0
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 21.66666666666667
$ws.Columns.Item(5).ColumnWidth = 16.33333333333333
$ws.Columns.Item(6).ColumnWidth = 23.83333333333333

# Move / fix the active selection
[void]$ws.Range("F4").Select()
